$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 240; existing rows 240:300 shift down to 241:301
$ws.Rows("240:240").Insert()

# Populate the newly inserted row 240 with the new record's data
$ws.Range("A240").Value = 3
$ws.Range("B240").Value = "Femacal de La Calera"
$ws.Range("C240").Value = "Coquimbo"
$ws.Range("D240").Value = 44711
$ws.Range("E240").Value = 5
$ws.Range("F240").Value = 100112013
$ws.Range("G240").Value = "Alcachofa"
$ws.Range("H240").Value = "Argentina(o)"
$ws.Range("I240").Value = "Primera"
$ws.Range("J240").Value = 50
$ws.Range("K240").Value = 15000
$ws.Range("L240").Value = 15000
$ws.Range("M240").Value = 15000
$ws.Range("N240").Value = "`$/caja 50 unidades"
$ws.Range("O240").Value = "Provincia de Quillota"
$ws.Range("P240").Value = 300
$ws.Range("Q240").Value = 50
$ws.Range("R240").Value = "Hortaliza"

# Match the date-formatted style used by the other D-column cells
$ws.Range("D240").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output "done"
